$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header E1: "me" -> "mase"
$ws.Range("E1").Value2 = "mase"

# Insert a new row above the current row 2 (shifts rows 2-8 down to 3-9)
$ws.Rows.Item(2).Insert()
$ws.Range("A2:E2").ClearFormats()

# Fill the new row 2 with the "NA" summary row
$ws.Range("A2").Value2 = "NA"
$ws.Range("B2").Value2 = 461.9178075874689
$ws.Range("C2").Value2 = 462.5952269423076
$ws.Range("D2").Value2 = 479.1689915266758
$ws.Range("E2").Value2 = 0.7131108469777713

# Update the recalculated "mase" column (E) for the shifted rows 3-9
$ws.Range("E3").Value2 = 0.7661727318958919
$ws.Range("E4").Value2 = 0.6365806040403128
$ws.Range("E5").Value2 = 0.7317147676616937
$ws.Range("E6").Value2 = 0.7437223847242481
$ws.Range("E7").Value2 = 0.6893679375550136
$ws.Range("E8").Value2 = 0.8104084128344408
$ws.Range("E9").Value2 = 0.7487955679452316
